# justification-type.docx: swap the East-Asian fallback font from
# "DejaVu Sans" to "Tahoma" on the styles that carry it directly
# (document defaults + the Normal/Heading paragraph styles), and make
# the List/Caption/Index styles carry an explicit complex-script
# ("cs") font of "DejaVu Sans" instead of inheriting it implicitly.

$d = $word.ActiveDocument

# --- East Asian fallback font: DejaVu Sans -> Tahoma -----------------
# (Normal and Heading both declare their own rFonts/eastAsia value.)
$normal = $d.Styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

$heading = $d.Styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

# --- Explicit complex-script (w:cs) font: DejaVu Sans -----------------
# List, Caption and Index previously had no direct rFonts override in
# their rPr; give them an explicit "cs" font matching the rest of the
# document instead of leaving it to fall back to the doc default.
$list = $d.Styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

$caption = $d.Styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

$index = $d.Styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
